$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 527.4286  # ALC!H2 (was 431.22223)
$ws.Cells.Item(2, 9).Value = 400  # ALC!I2 (was 277.8)
$ws.Cells.Item(2, 11).Value = 400  # ALC!K2 (was 277.8)
$ws.Cells.Item(2, 13).Value = -287  # ALC!M2 (was -164.8)
$ws.Cells.Item(20, 8).Value = 500  # ALC!H20 (was 5100)
$ws.Cells.Item(20, 9).Value = 500  # ALC!I20 (was 400)
$ws.Cells.Item(20, 10).Value = 0  # ALC!J20 (was 9800)
$ws.Cells.Item(20, 11).Value = 500  # ALC!K20 (was 400)
$ws.Cells.Item(20, 12).Value = 0  # ALC!L20 (was 9800)
$ws.Cells.Item(20, 13).Value = -270  # ALC!M20 (was -170)
$ws.Cells.Item(20, 14).Value = $null  # ALC!N20 (was -10260)
$ws.Cells.Item(35, 8).Value = 500  # ALC!H35 (was 5100)
$ws.Cells.Item(35, 9).Value = 500  # ALC!I35 (was 400)
$ws.Cells.Item(35, 10).Value = 0  # ALC!J35 (was 9800)
$ws.Cells.Item(35, 11).Value = 500  # ALC!K35 (was 400)
$ws.Cells.Item(35, 12).Value = 0  # ALC!L35 (was 9800)
$ws.Cells.Item(35, 13).Value = -121  # ALC!M35 (was -21)
$ws.Cells.Item(35, 14).Value = $null  # ALC!N35 (was -10558)
$ws.Cells.Item(40, 9).Value = 788  # ALC!I40 (was 781.6667)
$ws.Cells.Item(40, 10).Value = 1641.5834  # ALC!J40 (was 1722.6364)
$ws.Cells.Item(40, 11).Value = 788  # ALC!K40 (was 781.6667)
$ws.Cells.Item(40, 12).Value = 1641.5834  # ALC!L40 (was 1722.6364)
$ws.Cells.Item(40, 13).Value = -613  # ALC!M40 (was -606.6667)
$ws.Cells.Item(40, 14).Value = -1991.5834  # ALC!N40 (was -2072.6364)
$ws.Cells.Item(53, 8).Value = 4192.7  # ALC!H53 (was 3857)
$ws.Cells.Item(53, 9).Value = 356.66666  # ALC!I53 (was 392.5)
$ws.Cells.Item(53, 11).Value = 356.66666  # ALC!K53 (was 392.5)
$ws.Cells.Item(53, 13).Value = 280.33334  # ALC!M53 (was 244.5)
$ws.Cells.Item(64, 8).Value = 3399.9  # ALC!H64 (was 3061.4614)
$ws.Cells.Item(64, 9).Value = 3500  # ALC!I64 (was 2760)
$ws.Cells.Item(64, 10).Value = 3388.7778  # ALC!J64 (was 3249.875)
$ws.Cells.Item(64, 11).Value = 3500  # ALC!K64 (was 2760)
$ws.Cells.Item(64, 12).Value = 3388.7778  # ALC!L64 (was 3249.875)
$ws.Cells.Item(64, 13).Value = -3252  # ALC!M64 (was -2512)
$ws.Cells.Item(64, 14).Value = -3884.7778  # ALC!N64 (was -3745.875)
$ws.Cells.Item(67, 8).Value = 3399.9  # ALC!H67 (was 3061.4614)
$ws.Cells.Item(67, 9).Value = 3500  # ALC!I67 (was 2760)
$ws.Cells.Item(67, 10).Value = 3388.7778  # ALC!J67 (was 3249.875)
$ws.Cells.Item(67, 11).Value = 3500  # ALC!K67 (was 2760)
$ws.Cells.Item(67, 12).Value = 3388.7778  # ALC!L67 (was 3249.875)
$ws.Cells.Item(67, 13).Value = -2642  # ALC!M67 (was -1902)
$ws.Cells.Item(67, 14).Value = -5104.7778  # ALC!N67 (was -4965.875)
$ws.Cells.Item(70, 8).Value = 1500  # ALC!H70 (was 926.1579)
$ws.Cells.Item(70, 9).Value = 0  # ALC!I70 (was 842.4286)
$ws.Cells.Item(70, 10).Value = 1500  # ALC!J70 (was 975)
$ws.Cells.Item(70, 11).Value = 0  # ALC!K70 (was 2527.2858)
$ws.Cells.Item(70, 12).Value = 4500  # ALC!L70 (was 2925)
$ws.Cells.Item(70, 13).Value = $null  # ALC!M70 (was -2257.2858)
$ws.Cells.Item(70, 14).Value = -5040  # ALC!N70 (was -3465)
$ws.Cells.Item(73, 8).Value = 1500  # ALC!H73 (was 926.1579)
$ws.Cells.Item(73, 9).Value = 0  # ALC!I73 (was 842.4286)
$ws.Cells.Item(73, 10).Value = 1500  # ALC!J73 (was 975)
$ws.Cells.Item(73, 11).Value = 0  # ALC!K73 (was 2527.2858)
$ws.Cells.Item(73, 12).Value = 4500  # ALC!L73 (was 2925)
$ws.Cells.Item(73, 13).Value = $null  # ALC!M73 (was -1591.2858)
$ws.Cells.Item(73, 14).Value = -6372  # ALC!N73 (was -4797)
$ws.Cells.Item(92, 8).Value = 523.625  # ALC!H92 (was 616.0625)
$ws.Cells.Item(92, 9).Value = 521  # ALC!I92 (was 700.1818)
$ws.Cells.Item(92, 10).Value = 535  # ALC!J92 (was 431)
$ws.Cells.Item(92, 11).Value = 521  # ALC!K92 (was 700.1818)
$ws.Cells.Item(92, 12).Value = 535  # ALC!L92 (was 431)
$ws.Cells.Item(92, 13).Value = 727  # ALC!M92 (was 547.8182)
$ws.Cells.Item(92, 14).Value = -3031  # ALC!N92 (was -2927)
$ws.Cells.Item(116, 8).Value = 16670780  # ALC!H116 (was 15629047)
$ws.Cells.Item(116, 9).Value = 50000280  # ALC!I116 (was 83333816)
$ws.Cells.Item(116, 10).Value = 6029.6  # ALC!J116 (was 4868.923)
$ws.Cells.Item(116, 11).Value = 50000280  # ALC!K116 (was 83333816)
$ws.Cells.Item(116, 12).Value = 6029.6  # ALC!L116 (was 4868.923)
$ws.Cells.Item(116, 13).Value = -49996838  # ALC!M116 (was -83330374)
$ws.Cells.Item(116, 14).Value = -12913.6  # ALC!N116 (was -11752.923)
$ws.Cells.Item(129, 8).Value = 127488.76  # ALC!H129 (was 167613.2)
$ws.Cells.Item(129, 10).Value = 148072.73  # ALC!J129 (was 205186.25)
$ws.Cells.Item(129, 12).Value = 444218.1900000001  # ALC!L129 (was 615558.75)
$ws.Cells.Item(129, 14).Value = -454218.1900000001  # ALC!N129 (was -625558.75)

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(22, 8).Value = 1590.8334  # ARM!H22 (was 1257.625)
$ws.Cells.Item(22, 9).Value = 1590.8334  # ARM!I22 (was 1257.625)
$ws.Cells.Item(22, 11).Value = 1590.8334  # ARM!K22 (was 1257.625)
$ws.Cells.Item(22, 13).Value = -1291.8334  # ARM!M22 (was -958.625)
$ws.Cells.Item(32, 8).Value = 6343.9043  # ARM!H32 (was 5752.4487)
$ws.Cells.Item(32, 9).Value = 4429.0884  # ARM!I32 (was 4376.6377)
$ws.Cells.Item(32, 10).Value = 32385.4  # ARM!J32 (was 16300.333)
$ws.Cells.Item(32, 11).Value = 4429.0884  # ARM!K32 (was 4376.6377)
$ws.Cells.Item(32, 12).Value = 32385.4  # ARM!L32 (was 16300.333)
$ws.Cells.Item(32, 13).Value = -4142.0884  # ARM!M32 (was -4089.6377)
$ws.Cells.Item(32, 14).Value = -32959.4  # ARM!N32 (was -16874.333)
$ws.Cells.Item(63, 8).Value = 2084809.2  # ARM!H63 (was 1839797)
$ws.Cells.Item(63, 9).Value = 1634.5385  # ARM!I63 (was 1659.3125)
$ws.Cells.Item(63, 10).Value = 15625445  # ARM!J63 (was 31250000)
$ws.Cells.Item(63, 11).Value = 1634.5385  # ARM!K63 (was 1659.3125)
$ws.Cells.Item(63, 12).Value = 15625445  # ARM!L63 (was 31250000)
$ws.Cells.Item(63, 13).Value = -948.5385000000001  # ARM!M63 (was -973.3125)
$ws.Cells.Item(63, 14).Value = -15626817  # ARM!N63 (was -31251372)
$ws.Cells.Item(66, 8).Value = 2084809.2  # ARM!H66 (was 1839797)
$ws.Cells.Item(66, 9).Value = 1634.5385  # ARM!I66 (was 1659.3125)
$ws.Cells.Item(66, 10).Value = 15625445  # ARM!J66 (was 31250000)
$ws.Cells.Item(66, 11).Value = 8172.692500000001  # ARM!K66 (was 8296.5625)
$ws.Cells.Item(66, 12).Value = 78127225  # ARM!L66 (was 156250000)
$ws.Cells.Item(66, 13).Value = -4740.692500000001  # ARM!M66 (was -4864.5625)
$ws.Cells.Item(66, 14).Value = -78134089  # ARM!N66 (was -156256864)

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(15, 8).Value = 7800  # BSM!H15 (was 0)
$ws.Cells.Item(15, 10).Value = 7800  # BSM!J15 (was 0)
$ws.Cells.Item(15, 12).Value = 7800  # BSM!L15 (was 0)
$ws.Cells.Item(15, 14).Value = -8254  # BSM!N15 (was None)
$ws.Cells.Item(54, 8).Value = 9512.9  # BSM!H54 (was 8193.25)
$ws.Cells.Item(54, 9).Value = 2425.8  # BSM!I54 (was 2909.2)
$ws.Cells.Item(54, 10).Value = 16600  # BSM!J54 (was 17000)
$ws.Cells.Item(54, 11).Value = 2425.8  # BSM!K54 (was 2909.2)
$ws.Cells.Item(54, 12).Value = 16600  # BSM!L54 (was 17000)
$ws.Cells.Item(54, 13).Value = -1941.8  # BSM!M54 (was -2425.2)
$ws.Cells.Item(54, 14).Value = -17568  # BSM!N54 (was -17968)
$ws.Cells.Item(86, 8).Value = 1856.8206  # BSM!H86 (was 1859.3846)
$ws.Cells.Item(86, 9).Value = 1646.1923  # BSM!I86 (was 1650.0385)
$ws.Cells.Item(86, 11).Value = 1646.1923  # BSM!K86 (was 1650.0385)
$ws.Cells.Item(86, 13).Value = -523.1922999999999  # BSM!M86 (was -527.0385000000001)
$ws.Cells.Item(89, 8).Value = 1856.8206  # BSM!H89 (was 1859.3846)
$ws.Cells.Item(89, 9).Value = 1646.1923  # BSM!I89 (was 1650.0385)
$ws.Cells.Item(89, 11).Value = 8230.961499999999  # BSM!K89 (was 8250.192500000001)
$ws.Cells.Item(89, 13).Value = -2614.961499999999  # BSM!M89 (was -2634.192500000001)
$ws.Cells.Item(94, 8).Value = 692.0714  # BSM!H94 (was 705.0732)
$ws.Cells.Item(94, 9).Value = 454.5  # BSM!I94 (was 466.32)
$ws.Cells.Item(94, 11).Value = 454.5  # BSM!K94 (was 466.32)
$ws.Cells.Item(94, 13).Value = -3.5  # BSM!M94 (was -15.31999999999999)
$ws.Cells.Item(105, 8).Value = 2382807.2  # BSM!H105 (was 2633471)
$ws.Cells.Item(105, 9).Value = 1750  # BSM!I105 (was 1812.5)
$ws.Cells.Item(105, 10).Value = 3126887.5  # BSM!J105 (was 3335246.8)
$ws.Cells.Item(105, 11).Value = 1750  # BSM!K105 (was 1812.5)
$ws.Cells.Item(105, 12).Value = 3126887.5  # BSM!L105 (was 3335246.8)
$ws.Cells.Item(105, 13).Value = -3  # BSM!M105 (was -65.5)
$ws.Cells.Item(105, 14).Value = -3130381.5  # BSM!N105 (was -3338740.8)

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(141, 8).Value = 23556.525  # CRP!H141 (was 28040.523)
$ws.Cells.Item(141, 10).Value = 24122.219  # CRP!J141 (was 29042.55)
$ws.Cells.Item(141, 12).Value = 24122.219  # CRP!L141 (was 29042.55)
$ws.Cells.Item(141, 14).Value = -34482.219  # CRP!N141 (was -39402.55)

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(63, 8).Value = 4480.95  # CUL!H63 (was 4589.8887)
$ws.Cells.Item(63, 10).Value = 5101.1875  # CUL!J63 (was 5329.857)
$ws.Cells.Item(63, 12).Value = 15303.5625  # CUL!L63 (was 15989.571)
$ws.Cells.Item(63, 14).Value = -16801.5625  # CUL!N63 (was -17487.571)
$ws.Cells.Item(66, 8).Value = 4480.95  # CUL!H66 (was 4589.8887)
$ws.Cells.Item(66, 10).Value = 5101.1875  # CUL!J66 (was 5329.857)
$ws.Cells.Item(66, 12).Value = 45910.6875  # CUL!L66 (was 47968.713)
$ws.Cells.Item(66, 14).Value = -53398.6875  # CUL!N66 (was -55456.713)
$ws.Cells.Item(107, 8).Value = 5205  # CUL!H107 (was 4353.75)
$ws.Cells.Item(107, 9).Value = 6440.8125  # CUL!I107 (was 5443.421)
$ws.Cells.Item(107, 10).Value = 261.75  # CUL!J107 (was 213)
$ws.Cells.Item(107, 11).Value = 19322.4375  # CUL!K107 (was 16330.263)
$ws.Cells.Item(107, 12).Value = 785.25  # CUL!L107 (was 639)
$ws.Cells.Item(107, 13).Value = -17402.4375  # CUL!M107 (was -14410.263)
$ws.Cells.Item(107, 14).Value = -4625.25  # CUL!N107 (was -4479)
$ws.Cells.Item(117, 8).Value = 1356.2858  # CUL!H117 (was 1419.8462)
$ws.Cells.Item(117, 10).Value = 1355.9  # CUL!J117 (was 1447.6666)
$ws.Cells.Item(117, 12).Value = 4067.7  # CUL!L117 (was 4342.9998)
$ws.Cells.Item(117, 14).Value = -10951.7  # CUL!N117 (was -11226.9998)
$ws.Cells.Item(131, 8).Value = 703.22  # CUL!H131 (was 707.33)
$ws.Cells.Item(131, 9).Value = 457.14285  # CUL!I131 (was 483.33334)
$ws.Cells.Item(131, 10).Value = 721.74194  # CUL!J131 (was 721.6277)
$ws.Cells.Item(131, 11).Value = 1371.42855  # CUL!K131 (was 1450.00002)
$ws.Cells.Item(131, 12).Value = 2165.22582  # CUL!L131 (was 2164.8831)
$ws.Cells.Item(131, 13).Value = 3668.57145  # CUL!M131 (was 3589.99998)
$ws.Cells.Item(131, 14).Value = -12245.22582  # CUL!N131 (was -12244.8831)

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2022.3  # LTW!H22 (was 2040.6333)
$ws.Cells.Item(22, 9).Value = 2815.611  # LTW!I22 (was 2558.95)
$ws.Cells.Item(22, 10).Value = 832.3333  # LTW!J22 (was 1004)
$ws.Cells.Item(22, 11).Value = 2815.611  # LTW!K22 (was 2558.95)
$ws.Cells.Item(22, 12).Value = 832.3333  # LTW!L22 (was 1004)
$ws.Cells.Item(22, 13).Value = -2520.611  # LTW!M22 (was -2263.95)
$ws.Cells.Item(22, 14).Value = -1422.3333  # LTW!N22 (was -1594)
$ws.Cells.Item(27, 8).Value = 2022.3  # LTW!H27 (was 2040.6333)
$ws.Cells.Item(27, 9).Value = 2815.611  # LTW!I27 (was 2558.95)
$ws.Cells.Item(27, 10).Value = 832.3333  # LTW!J27 (was 1004)
$ws.Cells.Item(27, 11).Value = 2815.611  # LTW!K27 (was 2558.95)
$ws.Cells.Item(27, 12).Value = 832.3333  # LTW!L27 (was 1004)
$ws.Cells.Item(27, 13).Value = -2708.611  # LTW!M27 (was -2451.95)
$ws.Cells.Item(27, 14).Value = -1046.3333  # LTW!N27 (was -1218)

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(12, 8).Value = 80007  # WVR!H12 (was 100007)
$ws.Cells.Item(12, 10).Value = 80007  # WVR!J12 (was 100007)
$ws.Cells.Item(12, 12).Value = 80007  # WVR!L12 (was 100007)
$ws.Cells.Item(12, 14).Value = -80291  # WVR!N12 (was -100291)
$ws.Cells.Item(136, 8).Value = 25809518  # WVR!H136 (was 25180072)
$ws.Cells.Item(136, 9).Value = 33300006  # WVR!I136 (was 32259452)
$ws.Cells.Item(136, 11).Value = 99900018  # WVR!K136 (was 96778356)
$ws.Cells.Item(136, 13).Value = -99897468  # WVR!M136 (was -96775806)
